# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For every data row (rows 2-24), a new error value is inserted into
# column B and all the existing values in that row shift one column to
# the right (B->C, C->D, ... ), with anything that would spill past
# column K simply dropped. Row 1 (the Q0..Q9 header row) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to drop into column B for each row (2-24), keyed by row number.
$newValues = @{
    2  = 1.765134014560337
    3  = 7.384797829235454
    4  = -18.34129076817022
    5  = 7.600074697314557
    6  = 0.7916129955631771
    7  = -3.727363316492332
    8  = 0.376932102669816
    9  = 1.207578635508109
    10 = -0.9264868865757077
    11 = 0.3770345820039356
    12 = -0.4275923834192769
    13 = 0.324932645901923
    14 = -0.04071760298358112
    15 = 0.3721869518844864
    16 = -0.1524291232873974
    17 = -1.030518528898312
    18 = 0.4742145784871607
    19 = 0.3556547466179877
    20 = 0.3126006297022321
    21 = 0.3812981176718321
    22 = -0.716162849403934
    23 = 0.506656010950813
    24 = -0.343237405067616
}

$xlToLeft = -4159
$lastDataColumn = 30   # scan leftwards from a column past K to find the last occupied cell

for ($row = 2; $row -le 24; $row++) {

    # Last occupied value-column (1-based col index) in this row BEFORE the
    # edit - tells us where the right-shift needs to start from.
    $last = $ws.Cells.Item($row, $lastDataColumn).End($xlToLeft).Column
    if ($last -lt 2) { $last = 1 }

    # Shift existing values one column to the right, starting from the
    # rightmost occupied cell and working back down to column B (col 2),
    # so nothing is overwritten before it is read. Anything beyond column
    # K (col 11) falls off the edge and is discarded.
    for ($col = $last; $col -ge 2; $col--) {
        $val = $ws.Cells.Item($row, $col).Value2
        $destCol = $col + 1
        if ($destCol -le 11) {
            $ws.Cells.Item($row, $destCol).Value2 = $val
        }
    }

    # Drop the newly computed value into column B.
    $ws.Cells.Item($row, 2).Value2 = $newValues[$row]
}
